$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 6606
$ws.Range("I3").Value = 6897
$ws.Range("B4").Value = 1665
$ws.Range("G4").Value = 1447
$ws.Range("I4").Value = 1583
$ws.Range("I5").Value = 640
$ws.Range("I6").Value = 7927
$ws.Range("B7").Value = 23297
$ws.Range("G7").Value = 24672
$ws.Range("I7").Value = 23653
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 186
$ws.Range("I6").Value = 169
$ws.Range("I7").Value = 747
$ws.Range("I8").Value = 1415
$ws.Range("I10").Value = 170
$ws.Range("I15").Value = 276
$ws.Range("I19").Value = 669
$ws.Range("I23").Value = 232
$ws.Range("I26").Value = 31
$ws.Range("I29").Value = 1428
$ws.Range("I31").Value = 241
$ws.Range("I33").Value = 1060
$ws.Range("I41").Value = 102
$ws.Range("I42").Value = 862
$ws.Range("I44").Value = 176
$ws.Range("I46").Value = 52
$ws.Range("I50").Value = 119
$ws.Range("I51").Value = 283
$ws.Range("I52").Value = 520
$ws.Range("I53").Value = 259
$ws.Range("G54").Value = 335
$ws.Range("I54").Value = 477
$ws.Range("I55").Value = 269
$ws.Range("I60").Value = 133
$ws.Range("B63").Value = 369
$ws.Range("I63").Value = 79
$ws.Range("I65").Value = 542
$ws.Range("I67").Value = 911
$ws.Range("I70").Value = 38
$ws.Range("I72").Value = 94
$ws.Range("I78").Value = 319
$ws.Range("I79").Value = 670
$ws.Range("I84").Value = 210
$ws.Range("I85").Value = 1059
$ws.Range("I86").Value = 150
$ws.Range("I87").Value = 56
$ws.Range("I90").Value = 307
$ws.Range("I94").Value = 241
$ws.Range("I95").Value = 357
$ws.Range("I96").Value = 269
$ws.Range("I98").Value = 167
$ws.Range("B101").Value = 23297
$ws.Range("G101").Value = 24672
$ws.Range("I101").Value = 23653
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I3").Value = 404
$ws.Range("I6").Value = 274
$ws.Range("I7").Value = 1059
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I6").Value = 152
$ws.Range("I7").Value = 520
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 420
$ws.Range("I3").Value = 408
$ws.Range("I4").Value = 91
$ws.Range("I6").Value = 454
$ws.Range("I7").Value = 1415
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I2").Value = 58
$ws.Range("I7").Value = 259
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I3").Value = 231
$ws.Range("I7").Value = 747
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I6").Value = 103
$ws.Range("I7").Value = 269
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 338
$ws.Range("I6").Value = 276
$ws.Range("I7").Value = 911
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I2").Value = 73
$ws.Range("I7").Value = 241
$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I3").Value = 70
$ws.Range("I5").Value = 7
$ws.Range("I7").Value = 210
$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 182
$ws.Range("I7").Value = 542
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value = 125
$ws.Range("I7").Value = 357
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 239
$ws.Range("I3").Value = 390
$ws.Range("I7").Value = 1060
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("G4").Value = 22
$ws.Range("I6").Value = 231
$ws.Range("G7").Value = 335
$ws.Range("I7").Value = 477
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I3").Value = 495
$ws.Range("I7").Value = 1428
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 220
$ws.Range("I6").Value = 214
$ws.Range("I7").Value = 669
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("I3").Value = 52
$ws.Range("I7").Value = 176
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I6").Value = 48
$ws.Range("I7").Value = 169
$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I6").Value = 27
$ws.Range("I7").Value = 102
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I3").Value = 262
$ws.Range("I6").Value = 312
$ws.Range("I7").Value = 862
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("I3").Value = 34
$ws.Range("I7").Value = 170
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I2").Value = 72
$ws.Range("I7").Value = 319
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I3").Value = 86
$ws.Range("I7").Value = 269
$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("I2").Value = 14
$ws.Range("I3").Value = 18
$ws.Range("I7").Value = 52
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I2").Value = 64
$ws.Range("I3").Value = 82
$ws.Range("I7").Value = 232
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 193
$ws.Range("I7").Value = 670
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I2").Value = 45
$ws.Range("I7").Value = 241
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I3").Value = 65
$ws.Range("I7").Value = 276
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 106
$ws.Range("I7").Value = 167
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I2").Value = 35
$ws.Range("I7").Value = 119
$ws = $wb.Worksheets.Item('East Village')
$ws.Range("I6").Value = 20
$ws.Range("I7").Value = 31
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I2").Value = 68
$ws.Range("I3").Value = 60
$ws.Range("I7").Value = 186
$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("I3").Value = 11
$ws.Range("I7").Value = 38
$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I6").Value = 35
$ws.Range("I7").Value = 150
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I3").Value = 77
$ws.Range("I6").Value = 108
$ws.Range("I7").Value = 307
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I2").Value = 60
$ws.Range("I6").Value = 112
$ws.Range("I7").Value = 283
$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I2").Value = 46
$ws.Range("I7").Value = 133
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I6").Value = 45
$ws.Range("I7").Value = 94
$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("I6").Value = 32
$ws.Range("I7").Value = 56

Write-Host "Applied 146 cell updates for 2022-11-25 data"
